$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the NB row (was row 9) - deleting shifts SVM row up from 9 to 8
$ws.Rows("9:9").Delete()

# Extend the bold/border/center header style from C1:G1 (5 cells) onto the new H1:L1 headers
$ws.Range("C1:G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update header row (B1:L1)
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "One Year Base mean"
$ws.Range("D1").Value = "One Year Base std"
$ws.Range("E1").Value = "Two Year Base mean"
$ws.Range("F1").Value = "Two Year Base std"
$ws.Range("G1").Value = "Three Year Base mean"
$ws.Range("H1").Value = "Three Year Base std"
$ws.Range("I1").Value = "Five Year Base mean"
$ws.Range("J1").Value = "Five Year Base std"
$ws.Range("K1").Value = "Ten Year Base mean"
$ws.Range("L1").Value = "Ten Year Base std"

# Update algorithm names and data rows (rows 2-8)
$ws.Range("B2").Value = "LR"
$ws.Range("C2").Value = 0.8317426271571963
$ws.Range("D2").Value = 0.01715971997364202
$ws.Range("E2").Value = 0.8033368497016321
$ws.Range("F2").Value = 0.01306585200858986
$ws.Range("G2").Value = 0.7778297852612458
$ws.Range("H2").Value = 0.02064950253473409
$ws.Range("I2").Value = 0.7589321316852947
$ws.Range("J2").Value = 0.03640596425705438
$ws.Range("K2").Value = 0.7352025403434467
$ws.Range("L2").Value = 0.02607251922955917

$ws.Range("B3").Value = "LDA"
$ws.Range("C3").Value = 0.8325011319661169
$ws.Range("D3").Value = 0.01646979078870436
$ws.Range("E3").Value = 0.8085224903912558
$ws.Range("F3").Value = 0.01995440407397837
$ws.Range("G3").Value = 0.7861153457821016
$ws.Range("H3").Value = 0.01865946787760933
$ws.Range("I3").Value = 0.7625921613809259
$ws.Range("J3").Value = 0.04085995880965404
$ws.Range("K3").Value = 0.7377113586765425
$ws.Range("L3").Value = 0.03058972188679458

$ws.Range("B4").Value = "KNN"
$ws.Range("C4").Value = 0.8145307318468216
$ws.Range("D4").Value = 0.01952530995386
$ws.Range("E4").Value = 0.8052953430080002
$ws.Range("F4").Value = 0.01194983127508635
$ws.Range("G4").Value = 0.7988959717516018
$ws.Range("H4").Value = 0.02058333789712647
$ws.Range("I4").Value = 0.8000109637666627
$ws.Range("J4").Value = 0.02515356063336928
$ws.Range("K4").Value = 0.7801188808457724
$ws.Range("L4").Value = 0.0327547685636083

$ws.Range("B5").Value = "DTREE"
$ws.Range("C5").Value = 0.7773545747354136
$ws.Range("D5").Value = 0.02791660714363068
$ws.Range("E5").Value = 0.7725138683653101
$ws.Range("F5").Value = 0.01843029747341837
$ws.Range("G5").Value = 0.757278652895173
$ws.Range("H5").Value = 0.0149455206563262
$ws.Range("I5").Value = 0.7594109637817301
$ws.Range("J5").Value = 0.01908002487854356
$ws.Range("K5").Value = 0.7642742789976529
$ws.Range("L5").Value = 0.03982100080525339

$ws.Range("B6").Value = "RTREE"
$ws.Range("C6").Value = 0.8285672032604543
$ws.Range("D6").Value = 0.01249175156544549
$ws.Range("E6").Value = 0.8039313104964745
$ws.Range("F6").Value = 0.01629834113246902
$ws.Range("G6").Value = 0.7880214587958529
$ws.Range("H6").Value = 0.01208949883219567
$ws.Range("I6").Value = 0.7623704936737064
$ws.Range("J6").Value = 0.03389538437302407
$ws.Range("K6").Value = 0.7366169401383886
$ws.Range("L6").Value = 0.02939376754399653

$ws.Range("B7").Value = "XTREE"
$ws.Range("C7").Value = 0.8310201872203299
$ws.Range("D7").Value = 0.0180603276673216
$ws.Range("E7").Value = 0.8168284773799399
$ws.Range("F7").Value = 0.02633227974139423
$ws.Range("G7").Value = 0.8033637482405753
$ws.Range("H7").Value = 0.01691696951424706
$ws.Range("I7").Value = 0.7931524582033165
$ws.Range("J7").Value = 0.02516738916319707
$ws.Range("K7").Value = 0.7890413271233201
$ws.Range("L7").Value = 0.0218880873508531

$ws.Range("B8").Value = "SVM"
$ws.Range("C8").Value = 0.8238554171647789
$ws.Range("D8").Value = 0.01440361324830105
$ws.Range("E8").Value = 0.820511900741862
$ws.Range("F8").Value = 0.01628071281129193
$ws.Range("G8").Value = 0.8093517564401468
$ws.Range("H8").Value = 0.02041465758860431
$ws.Range("I8").Value = 0.8045937959344682
$ws.Range("J8").Value = 0.03361065626712906
$ws.Range("K8").Value = 0.7730943181097725
$ws.Range("L8").Value = 0.02980845192313035

